$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 61
$firstCol = 2   # B
$lastCol = 10   # J

$n = $lastRow - $firstRow + 1

# Read current values for each column across the row block, then write them back reversed.
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $colVals = @()
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $colVals += ,$ws.Cells.Item($r, $c).Value2
    }
    for ($i = 0; $i -lt $n; $i++) {
        $ws.Cells.Item($firstRow + $i, $c).Value2 = $colVals[$n - 1 - $i]
    }
}

Write-Output "done"
